# Rewrite the small "type demo" sheet into a year/sales table, per the
# commit "generate code by ai":
#   - header row becomes 年份 (Year) / 销售额 (Sales) instead of
#     整数/字符串/小数/日期 (int/string/decimal/date)
#   - row 2 keeps its existing numbers but the old text sample ("aa") is
#     replaced by a plain sales figure
#   - six rows of year/sales data are written (rows 2-7)
#   - the old demo columns C/D are cleared out (kept their formatting,
#     dropped their content)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "年份"

$ws.Range("B1").Value = "销售额"
$ws.Range("B1").Font.Name = "微软雅黑"

$ws.Range("C1").ClearContents()
$ws.Range("D1").ClearContents()

# --- Data rows --------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 20
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 3

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 40

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 50

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 100

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 60

# --- Leave the selection where the author left it --------------------
$null = $ws.Range("I5").Select()
